$wb = $excel.ActiveWorkbook

# The F column ("想去人数") values changed on both the "展览" and "全部类型"
# sheets (they mirror each other). Update the same cells on each.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 2081
    $ws.Range("F3").Value = 610
    $ws.Range("F4").Value = 1444
    $ws.Range("F5").Value = 6974
    $ws.Range("F7").Value = 109
}
